# Correct typos & update offloading fig
# Applies the 4 description-text corrections (added commas / trailing
# periods) to the "mxq_compile" parameter-reference sheet, and moves the
# viewport/selection the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- C2 : "model" parameter description -----------------------------------
# Add commas after each "Using backend=..." bullet line and a trailing
# period after the final line.
$modelDesc = @"
Model path or model instance. Model should be instance for the following cases:
Using backend="onnx" and a onnx model path,
Using backend="tvm" and a Keras model,
Using backend="tvm" and a PyTorch model,
Using backend="tf" and a fronzen TensorFlow PB graph.
"@
$ws.Range("C2").Value = $modelDesc

# --- C7 : "backend" parameter description ----------------------------------
# Add commas after each framework-mapping line and a trailing period.
$backendDesc = @"
Which framework to use to get the Mobilint IR.
It must be one of "onnx", "tf", and "tvm".
They correspond to deep learning frameworks as follows:
"onnx": ONNX,
"tf": TensorFlow,
"tvm": PyTorch, Keras, ONNX,
Defaults to "onnx".
"@
$ws.Range("C7").Value = $backendDesc

# --- C8 : "device" parameter description -----------------------------------
# Add a trailing period.
$deviceDesc = @"
Device to be used for compile and inerence. Either "cpu" or "gpu".
Defaults to "cpu".
"@
$ws.Range("C8").Value = $deviceDesc

# --- C10 : "quantize_percentile" parameter description ---------------------
# Add a trailing period.
$percentileDesc = @"
Percentile used for the quantization method "Percentile" and "MaxPercentile".
This should be between 0 and 1. (Ex. 0.999, 0.9999)
Defaults to 0.9999.
"@
$ws.Range("C10").Value = $percentileDesc

# --- Viewport / selection, matching the author's saved view ----------------
$ws.Activate()
try {
    $excel.ActiveWindow.ScrollRow = 7
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # view-state scrolling isn't always wired through COM; ignore if so
}
$ws.Range("C18").Select()
